# Fruta / hortaliza, semanal
# Insert two new weekly price rows (new row 3 and row 4), pushing the
# existing data rows (former rows 3-14) down to rows 5-16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows before the current row 3 (formatting of row 2 is
# copied down onto the new rows, including the date number format on
# column D).
$ws.Rows("3:4").Insert()

# ---- New row 3 ----
$ws.Cells.Item(3, 1).Value2 = 4
$ws.Cells.Item(3, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(3, 3).Value2 = "Los Lagos"
$ws.Cells.Item(3, 4).Value2 = 44537
$ws.Cells.Item(3, 5).Value2 = 10
$ws.Cells.Item(3, 6).Value2 = "Fruta"
$ws.Cells.Item(3, 7).Value2 = 100103
$ws.Cells.Item(3, 8).Value2 = "Frutos de hueso (carozo)"
$ws.Cells.Item(3, 9).Value2 = 100103003
$ws.Cells.Item(3, 10).Value2 = "Damasco"
$ws.Cells.Item(3, 11).Value2 = "Castle Brite"
$ws.Cells.Item(3, 12).Value2 = "Primera"
$ws.Cells.Item(3, 13).Value2 = 500
$ws.Cells.Item(3, 14).Value2 = 20000
$ws.Cells.Item(3, 15).Value2 = 22000
$ws.Cells.Item(3, 16).Value2 = 21000
$ws.Cells.Item(3, 17).Value2 = "`$/caja 18 kilos"
$ws.Cells.Item(3, 18).Value2 = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(3, 19).Value2 = 1167
$ws.Cells.Item(3, 20).Value2 = 18

# ---- New row 4 ----
$ws.Cells.Item(4, 1).Value2 = 4
$ws.Cells.Item(4, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(4, 3).Value2 = "Los Lagos"
$ws.Cells.Item(4, 4).Value2 = 44537
$ws.Cells.Item(4, 5).Value2 = 10
$ws.Cells.Item(4, 6).Value2 = "Fruta"
$ws.Cells.Item(4, 7).Value2 = 100103
$ws.Cells.Item(4, 8).Value2 = "Frutos de hueso (carozo)"
$ws.Cells.Item(4, 9).Value2 = 100103003
$ws.Cells.Item(4, 10).Value2 = "Damasco"
$ws.Cells.Item(4, 11).Value2 = "Castle Brite"
$ws.Cells.Item(4, 12).Value2 = "Segunda"
$ws.Cells.Item(4, 13).Value2 = 250
$ws.Cells.Item(4, 14).Value2 = 17000
$ws.Cells.Item(4, 15).Value2 = 17000
$ws.Cells.Item(4, 16).Value2 = 17000
$ws.Cells.Item(4, 17).Value2 = "`$/caja 18 kilos"
$ws.Cells.Item(4, 18).Value2 = "Región del Maule"
$ws.Cells.Item(4, 19).Value2 = 944
$ws.Cells.Item(4, 20).Value2 = 18

$dim = $ws.UsedRange.Address()
Write-Host ("Final UsedRange: " + $dim)
